$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update changed financial figures
$ws.Range("D2").Value = 470970
$ws.Range("E2").Value = 25725
$ws.Range("F2").Value = 25725
$ws.Range("G2").Value = 38163
$ws.Range("H2").Value = 29936
$ws.Range("I2").Value = 29936
$ws.Range("K2").Value = 410442
$ws.Range("L2").Value = 185603
$ws.Range("M2").Value = 224839
$ws.Range("N2").Value = 224839
$ws.Range("P2").Value = 21393
$ws.Range("Q2").Value = 23638
$ws.Range("R2").Value = -29834
$ws.Range("S2").Value = 9865
$ws.Range("T2").Value = 14296
$ws.Range("U2").Value = 9342
$ws.Range("V2").Value = 47010
$ws.Range("W2").Value = 5.46
$ws.Range("X2").Value = 6.36
$ws.Range("Y2").Value = 14.01
$ws.Range("Z2").Value = 7.75
$ws.Range("AA2").Value = 82.55
$ws.Range("AB2").Value = 960.67
$ws.Range("AC2").Value = 7385
$ws.Range("AD2").Value = 7.08
$ws.Range("AE2").Value = 55682
$ws.Range("AF2").Value = 0.9399999999999999
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 1.91
$ws.Range("AI2").Value = 13.5
$ws.Range("AJ2").Value = 405363347

# Row 3: update changed financial figures
$ws.Range("D3").Value = 495214
$ws.Range("E3").Value = 23543
$ws.Range("F3").Value = 23543
$ws.Range("G3").Value = 31003
$ws.Range("H3").Value = 26306
$ws.Range("I3").Value = 26306
$ws.Range("K3").Value = 459801
$ws.Range("L3").Value = 217761
$ws.Range("M3").Value = 242040
$ws.Range("N3").Value = 242040
$ws.Range("P3").Value = 21393
$ws.Range("Q3").Value = 33752
$ws.Range("R3").Value = -56138
$ws.Range("S3").Value = 9064
$ws.Range("T3").Value = 39146
$ws.Range("U3").Value = -5394
$ws.Range("V3").Value = 63174
$ws.Range("W3").Value = 4.75
$ws.Range("X3").Value = 5.31
$ws.Range("Y3").Value = 11.27
$ws.Range("Z3").Value = 6.05
$ws.Range("AA3").Value = 89.97
$ws.Range("AB3").Value = 1064.6
$ws.Range("AC3").Value = 6489
$ws.Range("AD3").Value = 8.109999999999999
$ws.Range("AE3").Value = 60369
$ws.Range("AF3").Value = 0.87
$ws.Range("AG3").Value = 1100
$ws.Range("AH3").Value = 2.09
$ws.Range("AI3").Value = 16.77
$ws.Range("AJ3").Value = 405363347

# Row 4: update changed financial figures
$ws.Range("D4").Value = 527129
$ws.Range("E4").Value = 24615
$ws.Range("F4").Value = 24615
$ws.Range("G4").Value = 34420
$ws.Range("H4").Value = 27546
$ws.Range("I4").Value = 27546
$ws.Range("K4").Value = 508893
$ws.Range("L4").Value = 243098
$ws.Range("M4").Value = 265794
$ws.Range("N4").Value = 265794
$ws.Range("P4").Value = 21393
$ws.Range("Q4").Value = 32759
$ws.Range("R4").Value = -23123
$ws.Range("S4").Value = 9454
$ws.Range("T4").Value = 15595
$ws.Range("U4").Value = 17163
$ws.Range("V4").Value = 80685
$ws.Range("W4").Value = 4.67
$ws.Range("X4").Value = 5.23
$ws.Range("Y4").Value = 10.85
$ws.Range("Z4").Value = 5.69
$ws.Range("AA4").Value = 91.45999999999999
$ws.Range("AB4").Value = 1177.08
$ws.Range("AC4").Value = 6795
$ws.Range("AD4").Value = 5.78
$ws.Range("AE4").Value = 66294
$ws.Range("AF4").Value = 0.59
$ws.Range("AG4").Value = 1100
$ws.Range("AH4").Value = 2.8
$ws.Range("AI4").Value = 16.01
$ws.Range("AJ4").Value = 405363347

# Row 5: update changed financial figures
$ws.Range("D5").Value = 535357
$ws.Range("E5").Value = 6622
$ws.Range("F5").Value = 6622
$ws.Range("G5").Value = 11401
$ws.Range("H5").Value = 9680
$ws.Range("I5").Value = 9680
$ws.Range("K5").Value = 522944
$ws.Range("L5").Value = 254333
$ws.Range("M5").Value = 268612
$ws.Range("N5").Value = 268612
$ws.Range("P5").Value = 21393
$ws.Range("Q5").Value = 25942
$ws.Range("R5").Value = -47946
$ws.Range("S5").Value = 7319
$ws.Range("T5").Value = 16895
$ws.Range("U5").Value = 9047
$ws.Range("V5").Value = 87530
$ws.Range("W5").Value = 1.24
$ws.Range("X5").Value = 1.81
$ws.Range("Y5").Value = 3.62
$ws.Range("Z5").Value = 1.88
$ws.Range("AA5").Value = 94.68000000000001
$ws.Range("AB5").Value = 1205.53
$ws.Range("AC5").Value = 2388
$ws.Range("AD5").Value = 14.03
$ws.Range("AE5").Value = 66997
$ws.Range("AF5").Value = 0.5
$ws.Range("AG5").Value = 800
$ws.Range("AH5").Value = 2.39
$ws.Range("AI5").Value = 33.13
$ws.Range("AJ5").Value = 405363347

# Row 6: update changed financial figures
$ws.Range("D6").Value = 541698
$ws.Range("E6").Value = 11575
$ws.Range("F6").Value = 11575
$ws.Range("G6").Value = 14686
$ws.Range("H6").Value = 11559
$ws.Range("I6").Value = 11559
$ws.Range("K6").Value = 517866
$ws.Range("L6").Value = 245431
$ws.Range("M6").Value = 272435
$ws.Range("N6").Value = 272435
$ws.Range("P6").Value = 21393
$ws.Range("Q6").Value = 44708
$ws.Range("R6").Value = -11554
$ws.Range("S6").Value = -25430
$ws.Range("T6").Value = 23763
$ws.Range("U6").Value = 20945
$ws.Range("V6").Value = 66837
$ws.Range("W6").Value = 2.14
$ws.Range("X6").Value = 2.13
$ws.Range("Y6").Value = 4.27
$ws.Range("Z6").Value = 2.22
$ws.Range("AA6").Value = 90.09
$ws.Range("AB6").Value = 1235.32
$ws.Range("AC6").Value = 2852
$ws.Range("AD6").Value = 11.82
$ws.Range("AE6").Value = 67950
$ws.Range("AF6").Value = 0.5
$ws.Range("AG6").Value = 900
$ws.Range("AH6").Value = 2.67
$ws.Range("AI6").Value = 31.22
$ws.Range("AJ6").Value = 405363347

# Row 7: update changed financial figures
$ws.Range("D7").Value = 568872
$ws.Range("E7").Value = 19690
$ws.Range("G7").Value = 27521
$ws.Range("H7").Value = 20421
$ws.Range("I7").Value = 20386
$ws.Range("K7").Value = 547156
$ws.Range("L7").Value = 256480
$ws.Range("M7").Value = 290676
$ws.Range("N7").Value = 289860
$ws.Range("P7").Value = 21391
$ws.Range("Q7").Value = 34564
$ws.Range("R7").Value = -18270
$ws.Range("S7").Value = -3617
$ws.Range("T7").Value = 21638
$ws.Range("U7").Value = 12182
$ws.Range("W7").Value = 3.46
$ws.Range("X7").Value = 3.59
$ws.Range("Y7").Value = 7.25
$ws.Range("Z7").Value = 3.84
$ws.Range("AA7").Value = 88.23999999999999
$ws.Range("AC7").Value = 5029
$ws.Range("AD7").Value = 8.23
$ws.Range("AE7").Value = 72297
$ws.Range("AF7").Value = 0.57
$ws.Range("AG7").Value = 961
$ws.Range("AH7").Value = 2.32
$ws.Range("AI7").Value = 19.11

# Row 8: update changed financial figures
$ws.Range("D8").Value = 611081
$ws.Range("E8").Value = 23464
$ws.Range("G8").Value = 30289
$ws.Range("H8").Value = 23167
$ws.Range("I8").Value = 23166
$ws.Range("K8").Value = 576988
$ws.Range("L8").Value = 269016
$ws.Range("M8").Value = 307972
$ws.Range("N8").Value = 306838
$ws.Range("P8").Value = 21391
$ws.Range("Q8").Value = 38791
$ws.Range("R8").Value = -27396
$ws.Range("S8").Value = -1084
$ws.Range("T8").Value = 20136
$ws.Range("U8").Value = 17864
$ws.Range("W8").Value = 3.84
$ws.Range("X8").Value = 3.79
$ws.Range("Y8").Value = 7.76
$ws.Range("Z8").Value = 4.12
$ws.Range("AA8").Value = 87.34999999999999
$ws.Range("AC8").Value = 5715
$ws.Range("AD8").Value = 7.16
$ws.Range("AE8").Value = 76531
$ws.Range("AF8").Value = 0.53
$ws.Range("AG8").Value = 1108
$ws.Range("AH8").Value = 2.71
$ws.Range("AI8").Value = 19.39

# Row 9: update changed financial figures
$ws.Range("D9").Value = 635502
$ws.Range("E9").Value = 26109
$ws.Range("G9").Value = 33954
$ws.Range("H9").Value = 25950
$ws.Range("I9").Value = 25950
$ws.Range("K9").Value = 606002
$ws.Range("L9").Value = 277764
$ws.Range("M9").Value = 328237
$ws.Range("N9").Value = 327682
$ws.Range("P9").Value = 21391
$ws.Range("Q9").Value = 42761
$ws.Range("R9").Value = -26609
$ws.Range("S9").Value = -1229
$ws.Range("T9").Value = 20310
$ws.Range("U9").Value = 23403
$ws.Range("W9").Value = 4.11
$ws.Range("X9").Value = 4.08
$ws.Range("Y9").Value = 8.18
$ws.Range("Z9").Value = 4.39
$ws.Range("AA9").Value = 84.62
$ws.Range("AC9").Value = 6402
$ws.Range("AD9").Value = 6.39
$ws.Range("AE9").Value = 81730
$ws.Range("AF9").Value = 0.5
$ws.Range("AG9").Value = 1179
$ws.Range("AH9").Value = 2.88
$ws.Range("AI9").Value = 18.42

# Columns J (net income, non-controlling) and O (equity, non-controlling)
# are removed for rows 2-5 (they did not exist for rows 6-9 already)
$ws.Range("J2:J5").ClearContents()
$ws.Range("O2:O5").ClearContents()
